# NOTCH.xlsx update: add "Lista de materiais" sheet (ADC.brd BOM sent to digikart)

$wb = $excel.ActiveWorkbook

# --- 1. Add the new worksheet at the end of the tab strip ---------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "Lista de materiais"

# --- 2. Column widths (engine stores ColumnWidth + 5/6 char, floored to
#        the nearest 1/6 char - compensate by subtracting 5/6 up front so
#        we land as close as possible to the real-Excel pixel-fit widths)
$off = 5.0 / 6.0
$ws.Columns.Item(1).ColumnWidth = 40 - $off
$ws.Range("B1:G1").EntireColumn.ColumnWidth = 4.5703125 - $off
$ws.Range("H1:K1").EntireColumn.ColumnWidth = 4 - $off
$ws.Columns.Item(12).ColumnWidth = 9.42578125 - $off
$ws.Columns.Item(13).ColumnWidth = 5.42578125 - $off
$ws.Columns.Item(14).ColumnWidth = 5.85546875 - $off
$ws.Columns.Item(15).ColumnWidth = 3.85546875 - $off
$ws.Columns.Item(16).ColumnWidth = 4.85546875 - $off
$ws.Columns.Item(17).ColumnWidth = 6.85546875 - $off
$ws.Range("R1:S1").EntireColumn.ColumnWidth = 7.5703125 - $off
$ws.Columns.Item(20).ColumnWidth = 15.5703125 - $off
$ws.Columns.Item(21).ColumnWidth = 16.5703125 - $off
$ws.Columns.Item(22).ColumnWidth = 8 - $off
$ws.Columns.Item(23).ColumnWidth = 15 - $off

# --- 3. Cell text/labels, written in the same order the original author
#        typed them (so new shared-string indices land in the same spots)
$ws.Range("A3").Value = "Ponte de Wheatstone "
$ws.Range("A4").Value = "Filtro 3Mhz"
$ws.Range("A5").Value = "Controle automático de ganho"
$ws.Range("A6").Value = "ADC"
$ws.Range("A7").Value = "DAC"
$ws.Range("A8").Value = "Avulsos"

$ws.Range("J2").Value = "22k"
$ws.Range("D2").Value = "4.7k"
$ws.Range("M2").Value = "4.7uf"
$ws.Range("N2").Value = "100nf"
$ws.Range("O2").Value = "1pf"
$ws.Range("P2").Value = "10pf"
$ws.Range("C2").Value = "3.3k"
$ws.Range("L2").Value = "PoT 100k"
$ws.Range("F2").Value = "6.8k"
$ws.Range("E2").Value = "5.6k"
$ws.Range("I2").Value = "12k"
$ws.Range("K2").Value = "27k"
$ws.Range("Q2").Value = "LT6323"
$ws.Range("R2").Value = "AD9288"
$ws.Range("S2").Value = "AD5445"
$ws.Range("V2").Value = "Jumper "
$ws.Range("T2").Value = "Barra pino(20x2)"
$ws.Range("U2").Value = "Born alimentação"
$ws.Range("W2").Value = "Barra pino (3x1)"
$ws.Range("B2").Value = "1.5k"

# G2 reuses the pre-existing "10k" shared string from Plan3
$ws.Range("G2").Value = "10k"

# O1 = "2.2" must stay TEXT (not the number 2.2). Force text storage via a
# temporary Text number-format, then ClearFormats so no stray style index
# is left behind (matches the un-styled target cell).
$ws.Range("O1").NumberFormat = "@"
$ws.Range("O1").Value = "2.2"
$ws.Range("O1").ClearFormats()

$ws.Range("K1").Value = "30k"
$ws.Range("J1").Value = "24k"
$ws.Range("I1").Value = "20k"
$ws.Range("H1").Value = "12.7k"
$ws.Range("F1").Value = "4.99k"
$ws.Range("E1").Value = "3.74k"
$ws.Range("D1").Value = "2k"

$ws.Range("V11").Value = "x'"

# --- 4. Quantity grid, rows 3-8 (centered, style matches existing "s=4")
#        Row 3 deliberately skips column H (matches the source workbook).
$ws.Range("B3:G3").HorizontalAlignment = -4108
$ws.Range("I3:W3").HorizontalAlignment = -4108
$ws.Range("G3").Value = 9

$ws.Range("B4:W4").HorizontalAlignment = -4108
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4
$ws.Range("N4").Value = 4
$ws.Range("P4").Value = 4
$ws.Range("Q4").Value = 1

$ws.Range("B5:W5").HorizontalAlignment = -4108
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 21
$ws.Range("Q5").Value = 3

$ws.Range("B6:W6").HorizontalAlignment = -4108
$ws.Range("R6").Value = 2

$ws.Range("B7:W7").HorizontalAlignment = -4108
$ws.Range("O7").Value = 1
$ws.Range("S7").Value = 1

$ws.Range("B8:W8").HorizontalAlignment = -4108
$ws.Range("T8").Value = 2
$ws.Range("U8").Value = 6
$ws.Range("V8").Value = 3
$ws.Range("W8").Value = 3

# --- 5. Totals row (plain, no special style) -----------------------------
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = 6
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 7
$ws.Range("G10").Value = 25
$ws.Range("H10").Value = 6
$ws.Range("I10").Value = 6
$ws.Range("J10").Value = 6
$ws.Range("K10").Value = 45
$ws.Range("M10").Value = 8
$ws.Range("N10").Value = 8
$ws.Range("O10").Value = 2
$ws.Range("P10").Value = 8

# --- 6. View state: select J14 (this also activates/selects this sheet) -
$ws.Activate()
$ws.Range("J14").Select()

# --- 7. Move workbook's active-tab pointer to this (4th, index 3) sheet -
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
